$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CCDeferredPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 22:25:34 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 22:26:37 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 19:29:42 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 19:30:54 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 21:01:47 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 21:02:53 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 20:57:24 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 20:58:43 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:48:00 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:48:50 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 22:15:42 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 22:16:58 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:51:28 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:52:18 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 22:20:52 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 22:21:56 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageDataCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Thu Aug 07 19:44:02 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Thu Aug 07 19:45:17 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:28:00 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:28:43 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 21:10:28 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 21:11:51 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 21:47:20 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 21:48:35 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:30:57 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:31:40 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:39:12 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:40:04 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 21:20:56 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 21:22:20 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:26:21 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:27:07 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 21:54:50 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 21:55:57 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 21:24:15 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 21:25:24 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 00:29:25 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 00:30:12 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 21:51:15 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 21:52:20 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 22:17:52 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Thu Aug 07 00:16:02 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 19:44:03 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 19:45:22 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 23:34:01 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 23:35:10 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 19:40:34 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 19:41:50 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 23:30:14 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 23:31:22 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 22:48:23 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 22:49:31 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 19:21:26 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 19:23:05 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 19:25:24 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 19:26:45 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 23:25:15 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 23:26:26 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 02:04:12 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 02:04:57 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 02:02:30 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 02:03:23 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 02:05:55 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 02:06:39 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 01:01:38 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 01:02:49 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 00:58:39 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 00:59:46 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 01:14:15 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 01:15:16 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 01:22:01 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 01:23:20 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 01:57:00 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 01:58:03 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 00:43:56 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 00:45:04 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 00:46:39 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 00:47:57 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Wed Aug 06 00:49:40 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Wed Aug 06 00:50:56 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Aug 05 02:07:29 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Aug 05 02:08:16 IST 2025"
